$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new route row ("ticker") to the Tabelle1 table -------------------
$tbl = $ws.ListObjects.Item("Tabelle1")
$newRow = $tbl.ListRows.Add()

$ws.Cells.Item(7, 1).Value = 3
$ws.Cells.Item(7, 2).Value = "/t"
$ws.Cells.Item(7, 3).Value = "ticker"
$ws.Cells.Item(7, 4).Value = "GET"
$ws.Cells.Item(7, 5).Value = "return basic info of a currency"

# --- Update the comment on the "price" route description (E6) --------------
$commentCell = $ws.Range("E6")
$newCommentText = "#GET price of BTC/USD, update every 10 seconds`nhttp://localhost:3000/p?from=BTC&to=USD`n#reponses:`n200, 500, 404"
$commentCell.Comment.Text($newCommentText)

# --- Leave the selection where the author left it off -----------------------
$ws.Range("G6").Select()
